$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old 4 data rows (rows 2-5), we'll rewrite 3 data rows (2-4)
$ws.Range("A2:D5").ClearContents()

# Deduplicated company data: Zara, Aleman, Fanta
$ws.Range("A2").Value = "Zara"
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 23
$ws.Range("D2").Value = "Empresa de ropa"

$ws.Range("A3").Value = "Aleman"
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = "Empresa de ropa"

$ws.Range("A4").Value = "Fanta"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "Empresa de alimentos"
